$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = 679025.8923849599
$ws.Range("F3").Value = 138914.46882130537
$ws.Range("G3").Value = 20.5

# Row 4
$ws.Range("E4").Value = 126789.67668903836
$ws.Range("F4").Value = 14516.597729148474
$ws.Range("G4").Value = 11.4

# Row 5
$ws.Range("E5").Value = 64519.75085621432
$ws.Range("F5").Value = 9845.141984185488
$ws.Range("G5").Value = 15.3

# Row 6
$ws.Range("E6").Value = 487716.46483970823
$ws.Range("F6").Value = 114552.72910797142
$ws.Range("G6").Value = 23.5

# Row 7
$ws.Range("E7").Value = 121243.75723927609
$ws.Range("F7").Value = 26125.173186287033
$ws.Range("G7").Value = 21.5

# Row 8
$ws.Range("E8").Value = 171097.71222410485
$ws.Range("F8").Value = 33849.03757102444
$ws.Range("G8").Value = 19.8

# Row 9
$ws.Range("E9").Value = 200942.87287758396
$ws.Range("F9").Value = 39263.165425997526
$ws.Range("G9").Value = 19.5

# Row 10
$ws.Range("E10").Value = 25238.457391499036
$ws.Range("F10").Value = 6052.406386566344
$ws.Range("G10").Value = 24.0

# Row 11
$ws.Range("E11").Value = 160372.4625551049
$ws.Range("F11").Value = 33595.9699673872
$ws.Range("G11").Value = 20.9

# Row 12
$ws.Range("E12").Value = 130.6300973915227
$ws.Range("F12").Value = 28.716284042823688
$ws.Range("G12").Value = 22.0
